$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.18
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 3.25
$ws.Range("J2").Value = 2.75
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 3.7
$ws.Range("N2").Value = 7.7
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 3.6
$ws.Range("Q2").Value = 1.75
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 1.37
$ws.Range("T2").Value = 2.85
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 2.25
$ws.Range("W2").Value = 8.75
$ws.Range("X2").Value = 11.75
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 16.5
$ws.Range("AB2").Value = 22
$ws.Range("AC2").Value = 7.7
$ws.Range("AD2").Value = 6.3
$ws.Range("AE2").Value = 11.5
$ws.Range("AF2").Value = 45
$ws.Range("AI2").Value = 19.5
$ws.Range("AK2").Value = 45
$ws.Range("AL2").Value = 26
$ws.Range("AM2").Value = 28
$ws.Range("AN2").Value = 4.25
$ws.Range("AO2").Value = 11.25
$ws.Range("AQ2").Value = 45
$ws.Range("AT2").Value = 2.85
$ws.Range("AU2").Value = 6.4
$ws.Range("AV2").Value = 50
$ws.Range("AW2").Value = 5.3
$ws.Range("AX2").Value = 17.5
$ws.Range("AY2").Value = 21
$ws.Range("AZ2").Value = 80
